# Generate Report for Handoff
# Adds a new row (row 9) for the file
# "dde28d1d-6d99-47bc-a6e3-891ef36efd0a.md" to the Overview / zh-cn / de-de
# sheets, and resizes each sheet's table/dimension accordingly.

$wb = $excel.ActiveWorkbook

$newFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67bca64a9c87545fdfe5882bd18bd99f09a1d888/e2e/dde28d1d-6d99-47bc-a6e3-891ef36efd0a.md"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A9").Value = "dde28d1d-6d99-47bc-a6e3-891ef36efd0a.md"
$wsOverview.Range("B9").Value = "e2e\dde28d1d-6d99-47bc-a6e3-891ef36efd0a.md"
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-18 06:42:36"
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), $newFileUrl, [Type]::Missing, [Type]::Missing, "e2e\dde28d1d-6d99-47bc-a6e3-891ef36efd0a.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G9"))

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A9").Value = "dde28d1d-6d99-47bc-a6e3-891ef36efd0a.md"
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "False"
$wsZhCn.Range("G9").Value = "dde28d1d-6d99-47bc-a6e3-891ef36efd0a.01354fe132c6ef1351d1be0b44fa9f34175c0e25.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-08-18 06:42:32"
$wsZhCn.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I9").Value = ""
$wsZhCn.Range("J9").Value = ""
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L9").Value = ""
$wsZhCn.Range("M9").Value = "True"
$wsZhCn.Range("N9").Value = ""
$wsZhCn.Range("O9").Value = "False"
$wsZhCn.Range("P9").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), $newFileUrl, [Type]::Missing, [Type]::Missing, "dde28d1d-6d99-47bc-a6e3-891ef36efd0a.md")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P9"))

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A9").Value = "dde28d1d-6d99-47bc-a6e3-891ef36efd0a.md"
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "False"
$wsDeDe.Range("G9").Value = "dde28d1d-6d99-47bc-a6e3-891ef36efd0a.01354fe132c6ef1351d1be0b44fa9f34175c0e25.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-08-18 06:42:36"
$wsDeDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I9").Value = ""
$wsDeDe.Range("J9").Value = ""
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L9").Value = ""
$wsDeDe.Range("M9").Value = "True"
$wsDeDe.Range("N9").Value = ""
$wsDeDe.Range("O9").Value = "False"
$wsDeDe.Range("P9").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), $newFileUrl, [Type]::Missing, [Type]::Missing, "dde28d1d-6d99-47bc-a6e3-891ef36efd0a.md")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P9"))
